# Generate Report for Archive
#
# The localization-status report advances the "Status" of every row from
# "Ready for handoff" to "In Translation". This touches:
#   - Overview sheet : columns E (zh-cn) and F (de-de), rows 2-4
#   - zh-cn sheet     : column C (Status), rows 2-4
#   - de-de sheet     : column C (Status), rows 2-4
# Because the new status text is shorter, the "Status" column on each of
# the three sheets is narrower (re-autofit after the text changed).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth = 13.4101845877511

# --- Overview sheet: Status columns are E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 2..4) {
    foreach ($col in @("E", "F")) {
        $cell = $wsOverview.Range("$col$row")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn / de-de sheets: Status column is C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in 2..4) {
        $cell = $ws.Range("C$row")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = $newWidth
}
